$d = $word.ActiveDocument

# The heading paragraph originally reads (across 4 runs):
#   "ตาราง" | " " | "… " | "Use case Description "
# and should become:
#   "ตาราง" | "ที่ 1" | " " | "Use case Description "
#
# Each Find/Replace below targets text that lives entirely within a single
# run, so the existing run boundaries (and their distinct formatting) are
# preserved.

# Run 3: "… " -> " "
$d.Content.Find.Execute("… ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, " ", 2)

# Run 2: the lone space run right after "ตาราง" -> "ที่ 1"
$d.Content.Find.Execute("ตาราง ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "ตารางที่ 1", 2)
